$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text figures (e.g. "43.575.65",
# "1.60") that must stay literal text rather than become parsed
# numbers. Briefly flip the cell to Text format while assigning the
# value, then restore General so the workbook formatting is untouched.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.575.65"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -6.55%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.592.44"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.82"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.37"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -4.42%  "

# Row 7
$ws.Range("E7").Value = "  -4.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -3.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.74"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -7.26%  "

# Row 11
$ws.Range("E11").Value = "  -3.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.81"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -4.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.997.01"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("E14").Value = "  +1.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.591.28"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.32"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -4.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.587.08"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -6.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.66"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -1.05%  "

# Row 20
$ws.Range("E20").Value = "  -4.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.31"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -5.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.84"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.56%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.61"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -5.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.58%  "

# Row 25
$ws.Range("E25").Value = "  -3.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.12"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +2.17%  "

# Row 27
$ws.Range("E27").Value = "  +0.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.86%  "

# Row 29 - swap with neighboring row
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.71"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.96%  "

# Row 30 - swap with neighboring row
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.16"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -6.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -5.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.61"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.20"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.55%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.99"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("E35").Value = "  -2.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0809"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -4.29%  "

# Row 37
$ws.Range("E37").Value = "  -4.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.45"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +5.52%  "

# Row 39
$ws.Range("E39").Value = "  -1.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.64"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +3.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42
$ws.Range("E42").Value = "  -5.56%  "

# Row 43
$ws.Range("E43").Value = "  -5.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.043.19"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").NumberFormat = "General"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.92"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -5.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -5.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.60"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +3.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.852.13"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.74"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.45%  "

# Row 51
$ws.Range("E51").Value = "  -5.19%  "
